$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Row 6: correct the F/G text (values were one row "off")
# ------------------------------------------------------------------
$ws.Range("F6").Value2 = "Change Performance timers"
$ws.Range("G6").Value2 = "Add in broadphase and narrowphase like in the tutorials"

# ------------------------------------------------------------------
# 2) Row 15: status TODO -> UNDERWAY (copy format from an existing
#    UNDERWAY-with-border cell, e.g. H3)
# ------------------------------------------------------------------
$ws.Range("H3").Copy()
$ws.Range("H15").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H15").Value2 = "UNDERWAY"

# ------------------------------------------------------------------
# 3) Row 16: update notes text; status TODO -> SATISFACTORY (new style:
#    same fill as SATISFACTORY-with-border H7, but without border)
# ------------------------------------------------------------------
$ws.Range("G16").Value2 = "Doesn't have any angular conditions"
$ws.Range("H7").Copy()
$ws.Range("H16").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H16").Borders.LineStyle = -4142 # xlNone - drop the border for this variant
$ws.Range("H16").Value2 = "SATISFACTORY"

# ------------------------------------------------------------------
# 4) Insert 5 fresh rows before the old row 17 (COMPOUND SHAPES),
#    pushing it down to row 22, to make room for new SOFT BODY tasks.
# ------------------------------------------------------------------
$ws.Range("17:21").Insert()

# Row 17: Generate a texturable mesh | UNDERWAY (no border)
$ws.Range("H14").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("H17").Value2 = "UNDERWAY"
$ws.Range("F17").Value2 = "Generate a texturable mesh"

# Row 18: Don't cull back faces | TODO (no border)
$ws.Range("H6").Copy()
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("H18").Value2 = "TODO"
$ws.Range("F18").Value2 = "Don't cull back faces"

# Row 19: Don't collide with itself | TODO (no border)
$ws.Range("H6").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("H19").Value2 = "TODO"
$ws.Range("F19").Value2 = "Don't collide with itself"

# Row 20: Make soft body class | TODO (plain)
$ws.Range("H2").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H20").Value2 = "TODO"
$ws.Range("F20").Value2 = "Make soft body class"

# Row 21: Make draggable | TODO (plain)
$ws.Range("H2").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("H21").Value2 = "TODO"
$ws.Range("F21").Value2 = "Make draggable"

# ------------------------------------------------------------------
# 5) Selection / active cell, per the target view state
# ------------------------------------------------------------------
[void]$ws.Range("L27").Select()

Write-Host "Done"
